# Reorder the rows of the "classFields" sheet so that the fields
# belonging to each class appear in the order produced by the
# (re-run) analysis tool, per the commit:
# "Implemented getting standard relationship between microservices
#  and started implementing MSM measure."
#
# The set of (Class Name, Field Name, Field Modifier, Field Type)
# rows is unchanged - only their order within each class block changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$rows = @(
    @("com.zatribune.spring.ecommerce.payments.db.DevBootstrap", "log", "private", "org.slf4j.Logger"),
    @("com.zatribune.spring.ecommerce.payments.db.DevBootstrap", "repository", "private", "com.zatribune.spring.ecommerce.payments.db.repository.CustomerRepository"),
    @("com.zatribune.spring.ecommerce.payments.listener.OrderListener", "log", "private", "org.slf4j.Logger"),
    @("com.zatribune.spring.ecommerce.payments.listener.OrderListener", "orderService", "private", "com.zatribune.spring.ecommerce.payments.service.OrderService"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer", "amountReserved", "private", "int"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer", "amountAvailable", "private", "int"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer", "id", "private", "java.lang.Long"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer", "name", "private", "java.lang.String"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer`$CustomerBuilder", "name", "private", "java.lang.String"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer`$CustomerBuilder", "id", "private", "java.lang.Long"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer`$CustomerBuilder", "amountAvailable", "private", "int"),
    @("com.zatribune.spring.ecommerce.payments.db.entities.Customer`$CustomerBuilder", "amountReserved", "private", "int"),
    @("com.zatribune.spring.ecommerce.payments.service.OrderServiceImpl", "repository", "private", "com.zatribune.spring.ecommerce.payments.db.repository.CustomerRepository"),
    @("com.zatribune.spring.ecommerce.payments.service.OrderServiceImpl", "template", "private", "org.springframework.kafka.core.KafkaTemplate"),
    @("com.zatribune.spring.ecommerce.payments.service.OrderServiceImpl", "log", "private", "org.slf4j.Logger"),
    @("com.zatribune.spring.ecommerce.payments.service.OrderServiceImpl", "SOURCE", "private", "domain.OrderSource")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
